$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# Supervisor Name
$ws.Range("G6").Value = "Prakruti Sinha"

# Supervisor Signature (initials) and sign-off date
$ws.Range("A27").Value = "P.S"
$ws.Range("D27").Value = (Get-Date -Year 2014 -Month 2 -Day 28 -Hour 0 -Minute 0 -Second 0).Date

$ws.Range("H33").Select() | Out-Null
